$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 410, shifting existing rows 410-456 down to 411-457.
$ws.Rows.Item(410).Insert()

# Populate the newly inserted row 410 with the new weekly price record.
$ws.Cells.Item(410, 1).Value = 8
$ws.Cells.Item(410, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(410, 3).Value = "Coquimbo"
$ws.Cells.Item(410, 4).Value = 44776
$ws.Cells.Item(410, 5).Value = 4
$ws.Cells.Item(410, 6).Value = 100112009
$ws.Cells.Item(410, 7).Value = "Acelga"
$ws.Cells.Item(410, 8).Value = "Sin especificar"
$ws.Cells.Item(410, 9).Value = "Segunda"
$ws.Cells.Item(410, 10).Value = 1360
$ws.Cells.Item(410, 11).Value = 600
$ws.Cells.Item(410, 12).Value = 650
$ws.Cells.Item(410, 13).Value = 625
$ws.Cells.Item(410, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(410, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(410, 16).Value = 312
$ws.Cells.Item(410, 17).Value = 2
$ws.Cells.Item(410, 18).Value = "Hortaliza"
